# Fixed Dev Gantt diagram
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Un-minimize the workbook window (workbookView minimized="1" removed).
$excel.WindowState = [Microsoft.Office.Interop.Excel.XlWindowState]::xlNormal
$wb.Windows.Item(1).WindowState = [Microsoft.Office.Interop.Excel.XlWindowState]::xlNormal

# 2. Update the task-4 ("Component unit testing") start/finish dates.
#    Moving F5/G5 forward causes the shared NETWORKDAYS formula in E5 to
#    recompute from "10g" to "25g" automatically.
$ws.Range("F5").Value2 = 42779.333333333336
$ws.Range("G5").Value2 = 42811.666666666664

# 3. Move the active selection on the sheet from B14 to G16.
$ws.Range("G16").Select()
